$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$range = $ws.Range("H2:H27")
$range.Value = "None"
$range.Select()
